# Update "想去人数" (want-to-go count) figures in the F column for the
# "展览" (sheet1) and "全部类型" (sheet4) worksheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 13655
$ws1.Range("F3").Value  = 87
$ws1.Range("F4").Value  = 124
$ws1.Range("F5").Value  = 538
$ws1.Range("F8").Value  = 1020
$ws1.Range("F9").Value  = 13871
$ws1.Range("F10").Value = 14683
$ws1.Range("F18").Value = 16
$ws1.Range("F19").Value = 55
$ws1.Range("F21").Value = 1138
$ws1.Range("F22").Value = 120
$ws1.Range("F24").Value = 5662
$ws1.Range("F26").Value = 1050
$ws1.Range("F27").Value = 5386
$ws1.Range("F29").Value = 42
$ws1.Range("F30").Value = 225

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 13655
$ws4.Range("F3").Value  = 87
$ws4.Range("F5").Value  = 124
$ws4.Range("F6").Value  = 538
$ws4.Range("F9").Value  = 1020
$ws4.Range("F10").Value = 13871
$ws4.Range("F11").Value = 14683
$ws4.Range("F19").Value = 16
$ws4.Range("F20").Value = 55
$ws4.Range("F22").Value = 1138
$ws4.Range("F23").Value = 120
$ws4.Range("F25").Value = 5662
$ws4.Range("F27").Value = 1050
$ws4.Range("F28").Value = 5386
$ws4.Range("F30").Value = 42
$ws4.Range("F31").Value = 225

$wb.Save()
